$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nuovi casi")
$ws.Cells.Item(450, 1).Value = 44348
$ws.Cells.Item(450, 3).Value = 1
$src = $ws.Range("D449")
$dst = $ws.Range("D450")
$src.Copy()
$dst.PasteSpecial(-4122)
Write-Host "done"
